$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "New" route's Path changed from "/home/new" to "/new"
$ws.Range("B3").Value = "/new"

# Reflect the new selection left on B3 after the edit
$ws.Range("B3").Select()

$wb.Save()
